$wb = $excel.ActiveWorkbook

# --- DatosCuenta (sheet1) ---
$ws1 = $wb.Worksheets.Item("DatosCuenta")
$ws1.Range("A2").Value = "SmokeName"
$ws1.Range("B2").Value = "SmokeLastName"
$ws1.Range("C2").Value = 20111100
$ws1.Range("D2").Value = 100

# --- DatosHogar (sheet2) ---
$ws2 = $wb.Worksheets.Item("DatosHogar")
$ws2.Range("A2").Value = 620

# --- DatosMotor (sheet3) ---
$ws3 = $wb.Worksheets.Item("DatosMotor")
$ws3.Range("A2").Value = "SMA001"
$ws3.Range("B2").Value = "ABC12SSMA001"
$ws3.Range("C2").Value = "ZAZ123SSMA001"

# --- DatosAP (sheet4) ---
$ws4 = $wb.Worksheets.Item("DatosAP")
$ws4.Range("A2").Value = 21200100

# Update selections on each sheet to match the saved view state
$ws1.Range("D3").Select() | Out-Null
$ws3.Range("C3").Select() | Out-Null

# Make DatosAP the active (selected) tab, with its selection on A3
$ws4.Activate()
$ws4.Range("A3").Select() | Out-Null
